# 7-2-1-1.xlsx metadata sheet edit
# - Re-style the column-A label cells (rows 2,3,4,6,7,8,9,10,12,13,14,16,17,19,20,21,23,24)
#   from the bordered-no-fill style to the bordered style used by B25 (s=1 in the
#   original workbook), matching the author's formatting pass over the label column.
# - Collapse the two blank "spacer" rows that separated the free-text footnote
#   (row 26) and the trailing comment cell (row 28) from their neighboring labeled
#   rows, so the sheet shrinks from 28 used rows down to 26.
# - Update the active selection / zoom to reflect the new layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Part 1: column A label-cell restyle -------------------------------
# B25 already carries the target format (fill + border + wrap, no bold) before
# any of the row surgery below, so use it as the format donor.
$ws.Range("B25").Copy()
$labelRows = @(2, 3, 4, 6, 7, 8, 9, 10, 12, 13, 14, 16, 17, 19, 20, 21, 23, 24)
foreach ($r in $labelRows) {
    $ws.Cells.Item($r, 1).PasteSpecial(-4122)   # xlPasteFormats
}

# --- Part 2: collapse rows 25-28 down to rows 25-26 ---------------------
# Old layout:
#   25: "Разбивка:" label        | blank input cell
#   26: (blank)                  | footnote text (shared string 33)
#   27: "8. Ссылки..." label     | blank input cell
#   28: (blank)                  | blank comment cell
# New layout:
#   25: "Разбивка:" label        | footnote text (shared string 33)
#   26: "8. Ссылки..." label     | blank comment cell

# Move row 26's value+format (footnote text) up into B25's placeholder, then
# delete the now-empty row 26.
$ws.Range("B26").Copy($ws.Range("B25"))
$ws.Rows("26:26").Delete()

# After the delete, old row 27 is row 26 and old row 28 is row 27. Move row 27's
# (the former row 28's) value+format into B26, then drop the spent row 27.
$ws.Range("B27").Copy($ws.Range("B26"))
$ws.Rows("27:27").Delete()

# --- Part 3: view state --------------------------------------------------
$ws.Range("A1:A26").Select()
$ws.Range("A26").Activate()
$excel.ActiveWindow.Zoom = 100
